# DeliveryNote.data.xlsx — add the "delivery" test-case rows (39-41) that
# carry an invalid dcid together with its related custid/dcno lookups.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39: invalid_dcid / 260
$ws.Range("A39").Value = "invalid_dcid"
$ws.Range("B39").Value = "260"

# Row 40: custid / 2566
$ws.Range("A40").Value = "custid"
$ws.Range("B40").Value = "2566"

# Row 41: dcno / 281
$ws.Range("A41").Value = "dcno"
$ws.Range("B41").Value = "281"

# Mirror the author's final cursor position: scrolled down so row 18 is at
# the top, with B41 (the last edited cell) selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 18
$win.ScrollColumn = 1

$ws.Range("B41").Select()
